$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Capture all existing values (rows 1-4, columns A-F) before we start shifting ----
$rowCount = 4
$colCount = 6
$data = @()
for ($r = 1; $r -le $rowCount; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $colCount; $c++) {
        $rowVals += ,$ws.Cells.Item($r, $c).Value2
    }
    $data += ,$rowVals
}

# ---- 2. Copy the header formatting (bold + border, from B1) onto A1 before we overwrite values ----
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ---- 3. Write the new header row: old B,C,D,E,F move one column left into A,B,C,D,E ----
#         and fix the MODEL_CONDITION -> MODELCONDITION text at the same time
$ws.Cells.Item(1,1).Value = $data[0][1]
$ws.Cells.Item(1,2).Value = $data[0][2]
$ws.Cells.Item(1,3).Value = $data[0][3]
$ws.Cells.Item(1,4).Value = "MODELCONDITION"
$ws.Cells.Item(1,5).Value = $data[0][5]

# ---- 4. Write the new data rows 2-4: old B,C,D,E,F move one column left into A,B,C,D,E ----
for ($r = 2; $r -le $rowCount; $r++) {
    $old = $data[$r - 1]
    $ws.Cells.Item($r,1).Value = $old[1]
    $ws.Cells.Item($r,2).Value = $old[2]
    $ws.Cells.Item($r,3).Value = $old[3]
    $ws.Cells.Item($r,4).Value = $old[4]
    $ws.Cells.Item($r,5).Value = $old[5]
}

# ---- 5. The old column A (rows 2-4) carried a bold/border style; the new column A data no longer should ----
$ws.Range("A2:A4").ClearFormats()

# ---- 6. Drop the now unused column F entirely ----
$ws.Range("F1:F4").Clear()

# ---- 7. Make sure the sheet's used range / dimension is tidy ----
$ws.Range("A1").Select() | Out-Null

Write-Host "Done"
